$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had a duplicated header row:
#   Row1: hours | mins | secs
#   Row2: hours | mins | secs   (duplicate header, to be removed)
#   Row3: 127   | 31   | 21
# Delete row 2 entirely so the data row (127/31/21) shifts up to become row 2.
$ws.Rows.Item(2).Delete()
